$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.31540000000003
$ws.Range("E5").Value = 13.22249999999999
$ws.Range("E9").Value = 13.94840000000002
$ws.Range("E11").Value = 13.52299999999999
$ws.Range("A21").Value = -21.28180000000001
$ws.Range("E21").Value = 12.6256
$ws.Range("A23").Value = -21.29530000000002
$ws.Range("A25").Value = -22.39300000000004
